$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "40-eridani-bc-astrometry"

# Insert new column B (band_eri_b) before current column B (ra_eri_b)
$ws.Columns("B").Insert()
# write "V" before the header so the shared-string table gets "V" at index 15
# and "band_eri_b" at index 16 (matches new-string append order)
$ws.Range("B2").Value = "V"
$ws.Range("B3").Value = "V"
$ws.Range("B4").Value = "V"
$ws.Range("B1").Value = "band_eri_b"

# Insert new column H (band_eri_c) before current column H (ra_eri_c, after first insert shift)
$ws.Columns("H").Insert()
$ws.Range("H1").Value = "band_eri_c"
$ws.Range("H2").Value = "I_c"
$ws.Range("H3").Value = "I_c"
$ws.Range("H4").Value = "I_c"

# Column widths (closest achievable to source widths 12.6640625 / 20.6640625)
$ws.Columns("B").ColumnWidth = 11.8333333333
$ws.Columns("H").ColumnWidth = 19.8333333333

# Selection
$ws.Range("H7").Select()
